# The author deleted row 13 (Genotype=LR2-36-01, trait=Size_maturity,
# treatment=PFOS, mean/sd blank) from Sheet1, causing every row below it
# to shift up by one (old row 14 -> new row 13, etc.) and the used range
# to shrink from A1:Q49 to A1:Q48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 13 and shift the rows below it up.
$ws.Rows.Item(13).Delete()

# Two other single-replicate rows (LR2-36-01/Fecundity/PFOS and
# LR2-36-01/Interval_brood/PFOS) had a blank "sd" value; after the shift
# they now live at rows 14 and 16 and their blank sd is filled in as 0.
$ws.Range("E14").Value = 0
$ws.Range("E16").Value = 0

# Leave the selection where Excel would land after such an edit.
$ws.Range("R24").Select()
